$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 data updates (re-run with fresh agenda/notes/date values)
$ws.Range("L3").Value = "Test agenda lagi"
$ws.Range("M3").Value = "notes"
$ws.Range("P3").Value = "2022-11-03"

# The "month" helper now reads the current month instead of 30 days ago
$ws.Range("Q3").Formula = "=TEXT(TODAY(),""mmmm"")"

# Move the viewport / selection like the re-run author did
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N2").Select()
